{"js": "// Ajuste na declara\u00e7\u00e3o do escopo\n// Applies the same textual changes as the reference OOXML diff for\n// \"1 - Declara\u00e7\u00e3o Escopo WigCred.docx\" using the Word JavaScript API.\n//\n// All edits are performed via Body.search() + Range.insertText(...,\n// Word.InsertLocation.Replace/Before) so that the existing run\n// formatting (the \"Hyperlink\" character style used throughout these\n// paragraphs) is preserved automatically.\n\nconst body = context.document.body;\n\n// 1) \"Wig Cred Cadastro...\" -> \"A Wig Cred Cadastro...\"\n//    (insert \"A \" immediately before the company name)\nlet startRange = body.search(\"Wig Cred Cadastro Cobran\u00e7a LTDA\", { matchCase: true });\nstartRange.load(\"text\");\nawait context.sync();\nif (startRange.items.length > 0) {\n  startRange.items[0].insertText(\"A \", Word.InsertLocation.before);\n  await context.sync();\n}\n\n// 2) \"...SP CEP 01004-010. A mais de 10 anos no mercado, a empresa hoje\n//    possui 6 funcion\u00e1rios...\"\n//    -> \"...SP CEP 01004-010, atua a mais de 10 anos no mercado, possui\n//    seis funcion\u00e1rios...\"\nlet addressRange = body.search(\n  \"01004-010. A mais de 10 anos no mercado, a empresa hoje possui 6 funcion\u00e1rios\",\n  { matchCase: true }\n);\naddressRange.load(\"text\");\nawait context.sync();\nif (addressRange.items.length > 0) {\n  addressRange.items[0].insertText(\n    \"01004-010, atua a mais de 10 anos no mercado, possui seis funcion\u00e1rios\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 3) \"...clientes e conta com o apoio de grandes parceiros como a\n//    Seekloc, SCPC e M\u00eddias Sociais.\" -> \"...clientes.\"\nlet partnersRange = body.search(\n  \"clientes e conta com o apoio de grandes parceiros como a Seekloc, SCPC e M\u00eddias Sociais.\",\n  { matchCase: true }\n);\npartnersRange.load(\"text\");\nawait context.sync();\nif (partnersRange.items.length > 0) {\n  partnersRange.items[0].insertText(\"clientes.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 4) \"A empresa Wig Cred \u00e9 especializada\" -> \"A empresa \u00e9 especializada\"\nlet companyRange = body.search(\"A empresa Wig Cred \u00e9 especializada\", { matchCase: true });\ncompanyRange.load(\"text\");\nawait context.sync();\nif (companyRange.items.length > 0) {\n  companyRange.items[0].insertText(\"A empresa \u00e9 especializada\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Ajuste na declara\u00e7\u00e3o do escopo\n# Applies the same textual changes as the reference OOXML diff for\n# \"1 - Declara\u00e7\u00e3o Escopo WigCred.docx\" using the Word COM object model.\n#\n# Each edit is done with Find (scoped to $d.Content, which Find collapses\n# to the matched span) followed by a direct Range.Text assignment, so the\n# surrounding run formatting (\"Hyperlink\" character style used throughout\n# these paragraphs) is preserved.\n\n$d = $word.ActiveDocument\n\nfunction Replace-FirstMatch($searchText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 0  # wdFindStop\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWildcards = $false\n    $rng.Find.Text = $searchText\n    $rng.Find.Execute() | Out-Null\n    if ($rng.Find.Found) {\n        $rng.Text = $replaceText\n    }\n}\n\n# 1) \"Wig Cred Cadastro...\" -> \"A Wig Cred Cadastro...\"\nReplace-FirstMatch \"Wig Cred Cadastro Cobran\u00e7a LTDA\" \"A Wig Cred Cadastro Cobran\u00e7a LTDA\"\n\n# 2) \"...SP CEP 01004-010. A mais de 10 anos no mercado, a empresa hoje\n#    possui 6 funcion\u00e1rios...\"\n#    -> \"...SP CEP 01004-010, atua a mais de 10 anos no mercado, possui\n#    seis funcion\u00e1rios...\"\nReplace-FirstMatch \"01004-010. A mais de 10 anos no mercado, a empresa hoje possui 6 funcion\u00e1rios\" \"01004-010, atua a mais de 10 anos no mercado, possui seis funcion\u00e1rios\"\n\n# 3) \"...clientes e conta com o apoio de grandes parceiros como a\n#    Seekloc, SCPC e M\u00eddias Sociais.\" -> \"...clientes.\"\nReplace-FirstMatch \"clientes e conta com o apoio de grandes parceiros como a Seekloc, SCPC e M\u00eddias Sociais.\" \"clientes.\"\n\n# 4) \"A empresa Wig Cred \u00e9 especializada\" -> \"A empresa \u00e9 especializada\"\nReplace-FirstMatch \"A empresa Wig Cred \u00e9 especializada\" \"A empresa \u00e9 especializada\"\n"}
